$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "Patient Auth." consent text in I2 ---
# The consent wording is replaced by a new version. The cell holds two runs:
# a bold quoted consent statement, followed by a normal explanatory sentence.
$boldText = '"My patient, and/or patient''s legal guardian if applicable, agrees to share their information with a study team, and to be contacted at the telephone number(s) and email they have provided, and via automated dialing, and/or artificial or pre-recorded voice, to schedule study appointments and keep them updated with important study-related information.  Their consent is not required as a condition of purchasing any property, goods, or services."'
$regularText = "  Options: Yes/No. To indicate whether the RP has the authorization from the patient or blank when not filled. Required Only if the file has Personal Information"

$cell = $ws.Range("I2")
$cell.Value = $boldText + $regularText

# First run (the quoted consent statement) is bold; keep the default size/font.
$cell.Characters(1, $boldText.Length).Font.Bold = $true
$cell.Characters(1, $boldText.Length).Font.Size = 12

# Second run (the explanatory sentence) stays regular weight.
$cell.Characters($boldText.Length + 1, $regularText.Length).Font.Bold = $false
$cell.Characters($boldText.Length + 1, $regularText.Length).Font.Size = 12

# --- Adjust the view: scroll position and active selection ---
$ws.Range("G5").Select()

# --- Column I gets its own (wider) width, split off from the F:H group ---
$ws.Columns("I").ColumnWidth = 15
